$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value (Excel serial number) for each
# data row (rows 2-98). The automatic update advanced this date by one day
# (45178 -> 45179) for every row.
$range = $ws.Range("C2:C98")
foreach ($cell in $range.Cells) {
    $cell.Value2 = [double]$cell.Value2 + 1
}
